# Update TPM-derived values in L1cam-Ephb2 sheet per new TPM computation
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 7.741029
$ws.Range("H2").Value = 23.223087
$ws.Range("I2").Value = 0.4930486933812723
$ws.Range("J2").Value = 0.4930486933812723
$ws.Range("M2").Value = 0.0006136666666666667
$ws.Range("N2").Value = 0.001841
$ws.Range("O2").Value = 0.000129696697123199
$ws.Range("P2").Value = 0.000129696697123199
$ws.Range("Q2").Value = 0.004750411463
$ws.Range("R2").Value = 0.042753703167
$ws.Range("S2").Value = 0.00006394678705245989
$ws.Range("T2").Value = 0.00006394678705245989
$ws.Range("G3").Value = 7.741029
$ws.Range("H3").Value = 23.223087
$ws.Range("I3").Value = 0.4930486933812723
$ws.Range("J3").Value = 0.4930486933812723
$ws.Range("O3").Value = 0.8077214410831794
$ws.Range("P3").Value = 0.8077214410831794
$ws.Range("Q3").Value = 29.584478847504
$ws.Range("R3").Value = 266.260309627536
$ws.Range("S3").Value = 0.3982460011420999
$ws.Range("T3").Value = 0.3982460011420999
$ws.Range("G4").Value = 7.741029
$ws.Range("H4").Value = 23.223087
$ws.Range("I4").Value = 0.4930486933812723
$ws.Range("J4").Value = 0.4930486933812723
$ws.Range("O4").Value = 0.1921488622196973
$ws.Range("P4").Value = 0.1921488622196973
$ws.Range("Q4").Value = 7.037851988041
$ws.Range("R4").Value = 63.340667892369
$ws.Range("S4").Value = 0.09473874545211988
$ws.Range("T4").Value = 0.09473874545211988
$ws.Range("I5").Value = 0.0194007766416684
$ws.Range("J5").Value = 0.0194007766416684
$ws.Range("M5").Value = 0.0006136666666666667
$ws.Range("N5").Value = 0.001841
$ws.Range("O5").Value = 0.000129696697123199
$ws.Range("P5").Value = 0.000129696697123199
$ws.Range("Q5").Value = 0.0001869220484444445
$ws.Range("R5").Value = 0.001682298436
$ws.Range("S5").Value = 0.000002516216652049301
$ws.Range("T5").Value = 0.000002516216652049301
$ws.Range("I6").Value = 0.0194007766416684
$ws.Range("J6").Value = 0.0194007766416684
$ws.Range("O6").Value = 0.8077214410831794
$ws.Range("P6").Value = 0.8077214410831794
$ws.Range("S6").Value = 0.01567042326714129
$ws.Range("T6").Value = 0.01567042326714129
$ws.Range("I7").Value = 0.0194007766416684
$ws.Range("J7").Value = 0.0194007766416684
$ws.Range("O7").Value = 0.1921488622196973
$ws.Range("P7").Value = 0.1921488622196973
$ws.Range("S7").Value = 0.003727837157875064
$ws.Range("T7").Value = 0.003727837157875064
$ws.Range("G8").Value = 7.654706000000001
$ws.Range("I8").Value = 0.4875505299770593
$ws.Range("J8").Value = 0.4875505299770593
$ws.Range("M8").Value = 0.0006136666666666667
$ws.Range("N8").Value = 0.001841
$ws.Range("O8").Value = 0.000129696697123199
$ws.Range("P8").Value = 0.000129696697123199
$ws.Range("Q8").Value = 0.004697437915333334
$ws.Range("R8").Value = 0.04227694123800001
$ws.Range("S8").Value = 0.00006323369341868984
$ws.Range("T8").Value = 0.00006323369341868984
$ws.Range("G9").Value = 7.654706000000001
$ws.Range("I9").Value = 0.4875505299770593
$ws.Range("J9").Value = 0.4875505299770593
$ws.Range("O9").Value = 0.8077214410831794
$ws.Range("P9").Value = 0.8077214410831794
$ws.Range("S9").Value = 0.3938050166739382
$ws.Range("T9").Value = 0.3938050166739382
$ws.Range("G10").Value = 7.654706000000001
$ws.Range("I10").Value = 0.4875505299770593
$ws.Range("J10").Value = 0.4875505299770593
$ws.Range("O10").Value = 0.1921488622196973
$ws.Range("P10").Value = 0.1921488622196973
$ws.Range("Q10").Value = 6.959370367940667
$ws.Range("R10").Value = 62.63433331146601
$ws.Range("S10").Value = 0.09368227960970238
$ws.Range("T10").Value = 0.09368227960970238
